# Update cached market-price / leve-profit figures across all job sheets
# (values refreshed by the scheduled market-data runner).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 2000
$ws.Range("K6").Value = 6000
$ws.Range("M6").Value = -5888
$ws.Range("H62").Value = 2503.88
$ws.Range("I62").Value = 2189.7856
$ws.Range("J62").Value = 2903.6365
$ws.Range("K62").Value = 2189.7856
$ws.Range("L62").Value = 2903.6365
$ws.Range("M62").Value = -1565.7856
$ws.Range("N62").Value = -4151.636500000001
$ws.Range("H65").Value = 2503.88
$ws.Range("I65").Value = 2189.7856
$ws.Range("J65").Value = 2903.6365
$ws.Range("K65").Value = 10948.928
$ws.Range("L65").Value = 14518.1825
$ws.Range("M65").Value = -7828.928
$ws.Range("N65").Value = -20758.1825
$ws.Range("H86").Value = 1291070.5
$ws.Range("I86").Value = 1939425.5
$ws.Range("J86").Value = 66399.89
$ws.Range("K86").Value = 1939425.5
$ws.Range("L86").Value = 66399.89
$ws.Range("M86").Value = -1938302.5
$ws.Range("N86").Value = -68645.89
$ws.Range("H89").Value = 1291070.5
$ws.Range("I89").Value = 1939425.5
$ws.Range("J89").Value = 66399.89
$ws.Range("K89").Value = 9697127.5
$ws.Range("L89").Value = 331999.45
$ws.Range("M89").Value = -9691511.5
$ws.Range("N89").Value = -343231.45
$ws.Range("H112").Value = 2717.07
$ws.Range("I112").Value = 1494
$ws.Range("J112").Value = 2785.0186
$ws.Range("K112").Value = 4482
$ws.Range("L112").Value = 8355.0558
$ws.Range("M112").Value = -3374
$ws.Range("N112").Value = -10571.0558

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 21230.645
$ws.Range("I32").Value = 25701.834
$ws.Range("J32").Value = 12288.267
$ws.Range("K32").Value = 25701.834
$ws.Range("L32").Value = 12288.267
$ws.Range("M32").Value = -25414.834
$ws.Range("N32").Value = -12862.267
$ws.Range("H45").Value = 2933.2144
$ws.Range("I45").Value = 1884.25
$ws.Range("J45").Value = 4331.8335
$ws.Range("K45").Value = 1884.25
$ws.Range("L45").Value = 4331.8335
$ws.Range("M45").Value = -1507.25
$ws.Range("N45").Value = -5085.8335
$ws.Range("H46").Value = 19082.467
$ws.Range("J46").Value = 19702.643
$ws.Range("L46").Value = 19702.643
$ws.Range("N46").Value = -20340.643
$ws.Range("H74").Value = 1222.9
$ws.Range("I74").Value = 1135.8
$ws.Range("J74").Value = 1484.2
$ws.Range("K74").Value = 1135.8
$ws.Range("L74").Value = 1484.2
$ws.Range("M74").Value = -261.8
$ws.Range("N74").Value = -3232.2
$ws.Range("H77").Value = 1222.9
$ws.Range("I77").Value = 1135.8
$ws.Range("J77").Value = 1484.2
$ws.Range("K77").Value = 5679
$ws.Range("L77").Value = 7421
$ws.Range("M77").Value = -1311
$ws.Range("N77").Value = -16157

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2448.7021
$ws.Range("I20").Value = 2019.4688
$ws.Range("J20").Value = 3364.4
$ws.Range("K20").Value = 2019.4688
$ws.Range("L20").Value = 3364.4
$ws.Range("M20").Value = -1772.4688
$ws.Range("N20").Value = -3858.4

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 457859.94
$ws.Range("I58").Value = 528308.9
$ws.Range("J58").Value = 404318.72
$ws.Range("K58").Value = 528308.9
$ws.Range("L58").Value = 404318.72
$ws.Range("M58").Value = -528105.9
$ws.Range("N58").Value = -404724.72
$ws.Range("H99").Value = 7714.1875
$ws.Range("I99").Value = 2889.25
$ws.Range("K99").Value = 2889.25
$ws.Range("M99").Value = -1391.25
$ws.Range("H126").Value = 7714.1875
$ws.Range("I126").Value = 2889.25
$ws.Range("K126").Value = 8667.75
$ws.Range("M126").Value = -6197.75
$ws.Range("H136").Value = 457859.94
$ws.Range("I136").Value = 528308.9
$ws.Range("J136").Value = 404318.72
$ws.Range("K136").Value = 1584926.7
$ws.Range("L136").Value = 1212956.16
$ws.Range("M136").Value = -1582376.7
$ws.Range("N136").Value = -1218056.16

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3732213.5
$ws.Range("I4").Value = 1286064.2
$ws.Range("J4").Value = 14006040
$ws.Range("K4").Value = 3858192.6
$ws.Range("L4").Value = 42018120
$ws.Range("M4").Value = -3858080.6
$ws.Range("N4").Value = -42018344
$ws.Range("H7").Value = 254.61539
$ws.Range("I7").Value = 207.625
$ws.Range("K7").Value = 622.875
$ws.Range("M7").Value = -510.875
$ws.Range("H41").Value = 416.33334
$ws.Range("I41").Value = 99
$ws.Range("J41").Value = 575
$ws.Range("K41").Value = 297
$ws.Range("L41").Value = 1725
$ws.Range("M41").Value = 41
$ws.Range("N41").Value = -2401
$ws.Range("H80").Value = 5893.5
$ws.Range("J80").Value = 5893.5
$ws.Range("L80").Value = 17680.5
$ws.Range("N80").Value = -19552.5
$ws.Range("H83").Value = 5893.5
$ws.Range("J83").Value = 5893.5
$ws.Range("L83").Value = 53041.5
$ws.Range("N83").Value = -62401.5
$ws.Range("H92").Value = 1496.2
$ws.Range("I92").Value = 1496.5
$ws.Range("K92").Value = 4489.5
$ws.Range("M92").Value = -3241.5
$ws.Range("H117").Value = 3414.2727
$ws.Range("J117").Value = 3810.5
$ws.Range("L117").Value = 11431.5
$ws.Range("N117").Value = -18315.5

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("N40").ClearContents()
$ws.Range("H45").Value = 28333.334
$ws.Range("J45").Value = 28333.334
$ws.Range("L45").Value = 28333.334
$ws.Range("N45").Value = -29451.334

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3072.963
$ws.Range("I16").Value = 1739.4
$ws.Range("K16").Value = 1739.4
$ws.Range("M16").Value = -1569.4
$ws.Range("H88").Value = 40189
$ws.Range("J88").Value = 40189
$ws.Range("L88").Value = 40189
$ws.Range("N88").Value = -41045
$ws.Range("H91").Value = 40189
$ws.Range("J91").Value = 40189
$ws.Range("L91").Value = 40189
$ws.Range("N91").Value = -43153
$ws.Range("H123").Value = 70108.664
$ws.Range("J123").Value = 70108.664
$ws.Range("L123").Value = 70108.664
$ws.Range("N123").Value = -79908.664

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3292.3333
$ws.Range("I126").Value = 2763.5
$ws.Range("J126").Value = 4350
$ws.Range("K126").Value = 8290.5
$ws.Range("L126").Value = 13050
$ws.Range("M126").Value = -5820.5
$ws.Range("N126").Value = -17990
$ws.Range("H129").Value = 85000
$ws.Range("J129").Value = 85000
$ws.Range("L129").Value = 85000
$ws.Range("N129").Value = -95000
